$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 27, shifting the existing rows 27-35 down to 29-37
$ws.Rows("27:28").Insert()

# Populate new row 27
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44855
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = "Otros"
$ws.Range("I27").Value = 100107002
$ws.Range("J27").Value = "Chirimoya"
$ws.Range("K27").Value = "Cultivar IV Región"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 110
$ws.Range("N27").Value = 24000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 24545
$ws.Range("Q27").Value = "$/bandeja 10 kilos"
$ws.Range("R27").Value = "Provincia de Limarí"
$ws.Range("S27").Value = 2454
$ws.Range("T27").Value = 10

# Populate new row 28
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44855
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = "Chirimoya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 22000
$ws.Range("O28").Value = 23000
$ws.Range("P28").Value = 22500
$ws.Range("Q28").Value = "$/bandeja 10 kilos"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 2250
$ws.Range("T28").Value = 10
